# Update "想去人数" (column F) values across the four sheets of the workbook.
# These numbers increased slightly between the two scrape/publish runs of the
# gh-pages data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F updates, keyed by row number.
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 3490
$ws.Cells.Item(5, 6).Value = 3490
$ws.Cells.Item(7, 6).Value = 5020
$ws.Cells.Item(8, 6).Value = 507
$ws.Cells.Item(10, 6).Value = 191
$ws.Cells.Item(14, 6).Value = 28
$ws.Cells.Item(15, 6).Value = 686
$ws.Cells.Item(16, 6).Value = 303
$ws.Cells.Item(19, 6).Value = 155
$ws.Cells.Item(22, 6).Value = 4864
$ws.Cells.Item(26, 6).Value = 5980
$ws.Cells.Item(30, 6).Value = 319
$ws.Cells.Item(31, 6).Value = 696
$ws.Cells.Item(34, 6).Value = 113
$ws.Cells.Item(36, 6).Value = 965
$ws.Cells.Item(37, 6).Value = 80
$ws.Cells.Item(40, 6).Value = 846
$ws.Cells.Item(41, 6).Value = 938

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 21

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 1106

# Sheet "全部类型" (All types) - aggregated listing with its own row numbering.
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 1106
$ws.Cells.Item(8, 6).Value = 3490
$ws.Cells.Item(9, 6).Value = 3490
$ws.Cells.Item(11, 6).Value = 5020
$ws.Cells.Item(12, 6).Value = 507
$ws.Cells.Item(14, 6).Value = 191
$ws.Cells.Item(17, 6).Value = 28
$ws.Cells.Item(18, 6).Value = 686
$ws.Cells.Item(19, 6).Value = 303
$ws.Cells.Item(23, 6).Value = 155
$ws.Cells.Item(26, 6).Value = 4864
$ws.Cells.Item(30, 6).Value = 5980
$ws.Cells.Item(34, 6).Value = 319
$ws.Cells.Item(35, 6).Value = 696
$ws.Cells.Item(38, 6).Value = 21
$ws.Cells.Item(39, 6).Value = 113
$ws.Cells.Item(41, 6).Value = 965
$ws.Cells.Item(42, 6).Value = 80
$ws.Cells.Item(45, 6).Value = 846
$ws.Cells.Item(46, 6).Value = 938
